# Sixth commit. Connecting logging to the project
#
# The underlying report data changed: the rows for the study profiles got
# reshuffled (profile "Физика" moved down to the last data row, "Медицина"
# moved up to the first data row, "Лингвистика" moved to the third row,
# while "Математика" stayed in place). Update the cell values in place to
# reflect the new row order; row/column styles are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now holds the "Медицина" profile data (previously in row 3)
$ws.Range("A2").Value = "Медицина"
$ws.Range("B2").Value = 4.300000190734863
$ws.Range("C2").Value = 3.0
$ws.Range("D2").Value = 3.0
$ws.Range("E2").Value = "Московский Государственный Медицинский Университет; Тамбовский Университет Медицины; Самарский Медицинский Институт; "

# Row 3: now holds the "Лингвистика" profile data (previously in row 5)
$ws.Range("A3").Value = "Лингвистика"
$ws.Range("B3").Value = 0.0
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 1.0
$ws.Range("E3").Value = "Воронежский Литературно-Переводческий Университет; "

# Row 4: "Математика" profile data is unchanged

# Row 5: now holds the "Физика" profile data (previously in row 2)
$ws.Range("A5").Value = "Физика"
$ws.Range("B5").Value = 4.5
$ws.Range("C5").Value = 8.0
$ws.Range("D5").Value = 2.0
$ws.Range("E5").Value = "Московский Выдуманный Университет; Московский Придуманный Институт; "
